# Feature_Tracker.xlsx update:
#  - Shortened two dice-related descriptions and cleared their "Never" status
#    (adding "Use Symbols preferably." + two new requestor attributions),
#    which also makes those two rows visible again under the existing filter.
#  - Row 8 ("Add Color") now has a completed version (1.7.7), which hides it
#    under the existing "blank = visible" filter on the Completed Version column.
#  - Added a couple of missing requestor attributions on two already-complete rows.
#  - Row 37 lost its (apparently incorrect / stale) Completed Version value.
#  - Two new feature requests were appended (rows 40 and 41).
#  - AutoFilter / selection / dimension naturally grow to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Add Color now shipped in 1.7.7 ---
$ws.Range("C8").Value = "1.7.7"

# --- Row 17: Add Genesys dice ---
$ws.Range("B17").Value = "The game Genesys uses dice. Add them in. Use Symbols preferably."
$ws.Range("C17").ClearContents()
$ws.Range("E17").Value = "forlasanto - Reddit /r/rpg"
$ws.Range("F17").Value = "IsaacNewtonsAndroid - /r/rpg"

# --- Row 18: Add Fantasy Flight Star Wars dice ---
$ws.Range("B18").Value = "The game Fantasy Flight Star Wars uses dice. Add them in. Use Symbols preferably."
$ws.Range("C18").ClearContents()
$ws.Range("E18").Value = "forlasanto - Reddit /r/rpg"
$ws.Range("F18").Value = "IsaacNewtonsAndroid - /r/rpg"

# --- Row 26: Make average more accurate - add requestor ---
$ws.Range("E26").Value = "michael - werbiskisfamily@gmail.com"

# --- Row 36: Dice with named faces - add requestors ---
$ws.Range("E36").Value = "forlasanto - Reddit /r/rpg"
$ws.Range("F36").Value = "IsaacNewtonsAndroid - /r/rpg"

# --- Row 37: Don't lose all rolls on uninstall - clear completed version ---
$ws.Range("C37").ClearContents()

# --- New row 40: Keep High/Low ---
$ws.Range("A40").Value = "Keep High/Low"
$ws.Range("B40").Value = "I have a variable number of dice rolled but only ever want to keep a set number"
$ws.Range("D40").Value = "michael - werbiskisfamily@gmail.com"

# --- New row 41: iOS support ---
$ws.Range("A41").Value = "iOS support"
$ws.Range("B41").Value = "I want to use the dice roller you made and put it on my iOS supported device"
$ws.Range("D41").Value = "All my iOS friends"

# --- Recompute row visibility to match the (already-applied) AutoFilter on
#     column C ("Completed Version"): filter keeps blanks visible, hides the
#     rest. Rows 8/17/18 flip visibility as a result of the edits above;
#     rows 40/41 are brand new and stay visible. ---
$ws.Rows.Item(8).Hidden = $true
$ws.Rows.Item(17).Hidden = $false
$ws.Rows.Item(18).Hidden = $false
$ws.Rows.Item(40).Hidden = $false
$ws.Rows.Item(41).Hidden = $false

# --- Re-apply the AutoFilter over the grown range so the persisted ref /
#     hidden _FilterDatabase name extend from E39 to E41 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:E41").AutoFilter(3, @(""))

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$41"
    }
}

# --- Match the author's final selection/active cell ---
$ws.Range("B44").Select()
